$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = "481 Plenty Rd, Preston VIC 3072"
$ws.Range("B6").Value = -37.73636
$ws.Range("C6").Value = 145.012539
$ws.Range("D6").Value = "Darebin (C)"
